# "fixed export and fixing maps"
#
# The sheet had been exported from a later/different source file that
# (a) carried a stray "1" sheet name instead of the municipality name,
# (b) added an extra "(census results)" note row, and
# (c) added two extra historical-census data columns (1989, 2002) that
#     don't belong on this particular map export - only the 2014 figure
#     should remain.
# This restores the sheet to the correct, trimmed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the sheet its proper (municipality) name instead of the generic "1".
$ws.Name = "დედოფლისწყარო"

# Drop the "(census results)" note row - row 2 in the original layout.
$ws.Rows.Item(2).Delete()

# Keep only the 2014 area figures; drop the 1989 and 2002 columns
# (originally columns B and C - the 2014 column shifts left into B).
$ws.Range("B:C").Delete()

# The 2014 header cell inherited the old "D-column" box border (thin
# left + medium right) from the column shift, but now that it's the
# only data column its box should look like the old left-most column
# header instead: no left rule, thin rule on the right.
$ws.Range("B4").Borders.Item(10).Weight = 2
$ws.Range("B4").Borders.Item(7).LineStyle = -4142

# Match the author's saved selection.
[void]$ws.Range("A2").Select()
